$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.764.57'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '2.232.41'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.85'
$ws.Range('E5').Value = '  +8.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.624'
$ws.Range('E6').Value = '  -2.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.64'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.18'
$ws.Range('E10').Value = '  +13.51%  '
$ws.Range('E11').Value = '  -2.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.51'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.98'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('D15').Value = '2.558.63'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.03'
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.856'
$ws.Range('E17').Value = '  -2.54%  '
$ws.Range('D18').Value = '2.234.18'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').Value = '41.752.17'
$ws.Range('E19').Value = '  -1.50%  '
$ws.Range('D20').Value = '0.0₃0969'
$ws.Range('E20').Value = '  -2.61%  '
$ws.Range('E21').Value = '  -1.10%  '
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.26'
$ws.Range('E23').Value = '  +15.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '233.08'
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').Value = '  +1.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.51'
$ws.Range('E27').Value = '  +5.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.23'
$ws.Range('E29').Value = '  +1.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.14'
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.76'
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('E33').Value = '  -1.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.48'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0721'
$ws.Range('E35').Value = '  -1.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.59'
$ws.Range('E36').Value = '  +15.61%  '
$ws.Range('E37').Value = '  -3.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.11'
$ws.Range('E38').Value = '  +11.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0284'
$ws.Range('E39').Value = '  +6.03%  '
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '70.29'
$ws.Range('E41').Value = '  +5.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.05'
$ws.Range('E42').Value = '  -1.76%  '
$ws.Range('E43').Value = '  +11.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.10'
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.72'
$ws.Range('E45').Value = '  +12.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.81'
$ws.Range('E47').Value = '  +8.04%  '
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('E50').Value = '  +6.51%  '
$ws.Range('D51').Value = '0.0₃0152'
$ws.Range('E51').Value = '  +13.39%  '
